$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells: force Text format so values like "29.248.83"
# or plain decimals like "325.93" are stored as literal text, matching
# the workbook's original inlineStr representation instead of being
# auto-coerced into numbers by Excel's input parser.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.248.83"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.10"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.93"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4636"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07876"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9886"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.81"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.909.78"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.062"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.731"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06989"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.25"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009970"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.04"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.250.70"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.148.13"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.099"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.85"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.38"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.002"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.32"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.889"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09354"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9017"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.255"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.323"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.207"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.185"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05780"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02084"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.712"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5708"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1786"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.701"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.90"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5355"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.182"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07024"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.847"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.566"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.96"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.047"

# Volume(1h) (column E) cells: plain text percentages, safe to assign directly.
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  -2.17%  "
